$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from AM1 into the new header cells AN1:AU1
$ws.Range("AM1").Copy()
$ws.Range("AN1:AU1").PasteSpecial(-4122)

# Set new header labels (row 1)
$ws.Cells.Item(1, 40).Value = "Diad1_Asym50"
$ws.Cells.Item(1, 41).Value = "Diad1_Asym70"
$ws.Cells.Item(1, 42).Value = "Diad1_Yuan2017_sym_factor"
$ws.Cells.Item(1, 43).Value = "Diad1_Remigi2021_BSF"
$ws.Cells.Item(1, 44).Value = "Diad2_Asym50"
$ws.Cells.Item(1, 45).Value = "Diad2_Asym70"
$ws.Cells.Item(1, 46).Value = "Diad2_Yuan2017_sym_factor"
$ws.Cells.Item(1, 47).Value = "Diad2_Remigi2021_BSF"

# Set new numeric data values (rows 2-17)
$ws.Cells.Item(2, 40).Value = 1.164133738602162
$ws.Cells.Item(2, 41).Value = 1.152173913043523
$ws.Cells.Item(2, 42).Value = 0.1922834017836951
$ws.Cells.Item(2, 43).Value = 0.0008571089448770777
$ws.Cells.Item(2, 44).Value = 1.100917431193286
$ws.Cells.Item(2, 45).Value = 1.100000000000032
$ws.Cells.Item(2, 46).Value = 0.1027222787810408
$ws.Cells.Item(2, 47).Value = 0.0004185099443979589

$ws.Cells.Item(3, 40).Value = 1.181008902076529
$ws.Cells.Item(3, 41).Value = 1.146551724137428
$ws.Cells.Item(3, 42).Value = 0.2139112653169737
$ws.Cells.Item(3, 43).Value = 0.0009137605674561974
$ws.Cells.Item(3, 44).Value = 1.081081081080656
$ws.Cells.Item(3, 45).Value = 1.081775700934434
$ws.Cells.Item(3, 46).Value = 0.08279008964124143
$ws.Cells.Item(3, 47).Value = 0.0004555797481215617

$ws.Cells.Item(4, 40).Value = 1.198198198198423
$ws.Cells.Item(4, 41).Value = 1.185185185185261
$ws.Cells.Item(4, 42).Value = 0.2350882955196955
$ws.Cells.Item(4, 43).Value = 0.0006390089314212232
$ws.Cells.Item(4, 44).Value = 1.06547619047562
$ws.Cells.Item(4, 45).Value = 1.067129629629676
$ws.Cells.Item(4, 46).Value = 0.06721912718399045
$ws.Cells.Item(4, 47).Value = 0.0003130828167376919

$ws.Cells.Item(5, 40).Value = 1.195035460993019
$ws.Cells.Item(5, 41).Value = 1.20111731843574
$ws.Cells.Item(5, 42).Value = 0.1829373502616893
$ws.Cells.Item(5, 43).Value = 0.02019587356651059
$ws.Cells.Item(5, 44).Value = 1.034482758620562
$ws.Cells.Item(5, 45).Value = 1.03888888888852
$ws.Cells.Item(5, 46).Value = 0.02574252342936353
$ws.Cells.Item(5, 47).Value = 0.009673999294495873

$ws.Cells.Item(6, 40).Value = 1.252631578946804
$ws.Cells.Item(6, 41).Value = 1.266846361185536
$ws.Cells.Item(6, 42).Value = 0.2436567203689979
$ws.Cells.Item(6, 43).Value = 0.001527742367305308
$ws.Cells.Item(6, 44).Value = 1.095057034220433
$ws.Cells.Item(6, 45).Value = 1.122507122507193
$ws.Cells.Item(6, 46).Value = 0.08687053902164894
$ws.Cells.Item(6, 47).Value = 0.0008885840855582906

$ws.Cells.Item(7, 40).Value = 1.201954397394592
$ws.Cells.Item(7, 41).Value = 1.201511335012676
$ws.Cells.Item(7, 42).Value = 0.2020226264871046
$ws.Cells.Item(7, 43).Value = 0.002408408218673916
$ws.Cells.Item(7, 44).Value = 1.06716417910465
$ws.Cells.Item(7, 45).Value = 1.088642659279925
$ws.Cells.Item(7, 46).Value = 0.06234023522498912
$ws.Cells.Item(7, 47).Value = 0.001366294128249545

$ws.Cells.Item(8, 40).Value = 1.003597122302353
$ws.Cells.Item(8, 41).Value = 1.002645502645645
$ws.Cells.Item(8, 42).Value = 0.003388154343009955
$ws.Cells.Item(8, 43).Value = 0.004003132965718732
$ws.Cells.Item(8, 44).Value = 1.251937984496013
$ws.Cells.Item(8, 45).Value = 1.28398791540778
$ws.Cells.Item(8, 46).Value = 0.2162185216442824
$ws.Cells.Item(8, 47).Value = 0.002122158418979997

$ws.Cells.Item(9, 40).Value = 1.017595307917657
$ws.Cells.Item(9, 41).Value = 1.027310924369817
$ws.Cells.Item(9, 42).Value = 0.02052391063450729
$ws.Cells.Item(9, 43).Value = 0.003601925827545564
$ws.Cells.Item(9, 44).Value = 1.221476510067057
$ws.Cells.Item(9, 45).Value = 1.22137404580167
$ws.Cells.Item(9, 46).Value = 0.227293600785559
$ws.Cells.Item(9, 47).Value = 0.001836225772890855

$ws.Cells.Item(10, 40).Value = 1.094674556212968
$ws.Cells.Item(10, 41).Value = 1.095541401273929
$ws.Cells.Item(10, 42).Value = 0.1123638294330259
$ws.Cells.Item(10, 43).Value = 0.0006100426740246082
$ws.Cells.Item(10, 44).Value = 1.197411003235526
$ws.Cells.Item(10, 45).Value = 1.190123456789631
$ws.Cells.Item(10, 46).Value = 0.2029290416787086
$ws.Cells.Item(10, 47).Value = 0.0003006820035996373

$ws.Cells.Item(11, 40).Value = 1.025568181818422
$ws.Cells.Item(11, 41).Value = 1.053388090348928
$ws.Cells.Item(11, 42).Value = 0.03077567239018532
$ws.Cells.Item(11, 43).Value = 0.0007172586770451544
$ws.Cells.Item(11, 44).Value = 1.203225806451472
$ws.Cells.Item(11, 45).Value = 1.193627450980329
$ws.Cells.Item(11, 46).Value = 0.2120469227446246
$ws.Cells.Item(11, 47).Value = 0.0003537435869755437

$ws.Cells.Item(12, 40).Value = 1.026086956521594
$ws.Cells.Item(12, 41).Value = 1.008196721311367
$ws.Cells.Item(12, 42).Value = 0.03089089429699293
$ws.Cells.Item(12, 43).Value = 0.001907789716365235
$ws.Cells.Item(12, 44).Value = 1.223333333333295
$ws.Cells.Item(12, 45).Value = 1.218592964823898
$ws.Cells.Item(12, 46).Value = 0.2310084308933085
$ws.Cells.Item(12, 47).Value = 0.0009513085296379747

$ws.Cells.Item(13, 40).Value = 1.022922636102848
$ws.Cells.Item(13, 41).Value = 1.010183299389007
$ws.Cells.Item(13, 42).Value = 0.02746718630240945
$ws.Cells.Item(13, 43).Value = 0.002555081465427714
$ws.Cells.Item(13, 44).Value = 1.217821782178126
$ws.Cells.Item(13, 45).Value = 1.212499999999719
$ws.Cells.Item(13, 46).Value = 0.2244890929404168
$ws.Cells.Item(13, 47).Value = 0.001258444643877306

$ws.Cells.Item(14, 40).Value = 1.025714285714142
$ws.Cells.Item(14, 41).Value = 1.033482142857151
$ws.Cells.Item(14, 42).Value = 0.02792417226723898
$ws.Cells.Item(14, 43).Value = 0.003273157864424728
$ws.Cells.Item(14, 44).Value = 1.038732394366109
$ws.Cells.Item(14, 45).Value = 1.035989717223308
$ws.Cells.Item(14, 46).Value = 0.03785469975302655
$ws.Cells.Item(14, 47).Value = 0.001747704912994523

$ws.Cells.Item(15, 40).Value = 1.101744186046138
$ws.Cells.Item(15, 41).Value = 1.068131868131711
$ws.Cells.Item(15, 42).Value = 0.1156179298024332
$ws.Cells.Item(15, 43).Value = 0.001701159774999307
$ws.Cells.Item(15, 44).Value = 1.03092783505186
$ws.Cells.Item(15, 45).Value = 1.047619047619528
$ws.Cells.Item(15, 46).Value = 0.03091937550983918
$ws.Cells.Item(15, 47).Value = 0.000919340054609549

$ws.Cells.Item(16, 40).Value = 1.016574585635515
$ws.Cells.Item(16, 41).Value = 1.008528784648363
$ws.Cells.Item(16, 42).Value = 0.01856866585007859
$ws.Cells.Item(16, 43).Value = 0.002757223678063495
$ws.Cells.Item(16, 44).Value = 1.03496503496491
$ws.Cells.Item(16, 45).Value = 1.043701799485593
$ws.Cells.Item(16, 46).Value = 0.03405848687296949
$ws.Cells.Item(16, 47).Value = 0.00152434322579094

$ws.Cells.Item(17, 40).Value = 1.007462686567234
$ws.Cells.Item(17, 41).Value = 1.070610687022875
$ws.Cells.Item(17, 42).Value = 0.01000199511166949
$ws.Cells.Item(17, 43).Value = 0.002807459033431672
$ws.Cells.Item(17, 44).Value = 1.195286195286329
$ws.Cells.Item(17, 45).Value = 1.20749999999967
$ws.Cells.Item(17, 46).Value = 0.2050119700710883
$ws.Cells.Item(17, 47).Value = 0.001203286975286477

